$wb = $excel.ActiveWorkbook

# --- Sheet "weeknr 48" (sheet1) -------------------------------------------
$ws1 = $wb.Worksheets.Item("weeknr 48")

# Rows 18, 20, 21, 22 had their Id (E) and duration (G) values removed,
# leaving only the formatted, empty cells behind.
$ws1.Range("E18").ClearContents()
$ws1.Range("E20:E22").ClearContents()
$ws1.Range("G18").ClearContents()
$ws1.Range("G20:G25").ClearContents()

# Row 19 was fully cleared (no data, no per-cell styling left at all).
$ws1.Rows("19:19").Clear()

# Selection on this sheet moved to G7 and it is no longer the active tab.
$ws1.Range("G7").Select()

# --- Sheet "week 49" (sheet2) ----------------------------------------------
$ws2 = $wb.Worksheets.Item("week 49")

# Fix the project name, which had been left showing an old value.
$ws2.Range("B2").Value = "PyramidPanic"

# Thursday's log entry: day label, date, description and duration formula.
$ws2.Range("A7").Value = "Donderdag"
$ws2.Range("B7").Value = 41612
$ws2.Range("F7").Value = "StartScene in PyramaidPanic gezet"
$ws2.Range("G7").Formula = "=D7-C7"

# A second entry was added on row 8.
$ws2.Range("C8").Value = 0.39583333333333331
$ws2.Range("D8").Value = 0.41319444444444442
$ws2.Range("E8").Value = 2
$ws2.Range("F8").Value = "De rest van de Scene's gemaakt"

# This sheet becomes the active / selected tab, with G8 selected.
$ws2.Activate()
$ws2.Range("G8").Select()
